$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, shifting the existing rows (17-159) down to (18-160)
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the "BETTER DRINKS" client record
$ws.Range("A17").Value = 174
$ws.Range("B17").Value = "BETTER DRINKS"
$ws.Range("V17").Value = 44700
$ws.Range("W17").Value = 0
$ws.Range("X17").Value = 0
$ws.Range("Y17").Value = 0
$ws.Range("AA17").Value = 0
$ws.Range("AB17").Value = $false

# Update the "Clientes" defined name to reflect the new last row (159 -> 160)
$n = $wb.Names.Item(1)
$n.RefersTo = "='Clientes'!`$A`$1:`$AE`$160"
